$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36, shifting existing rows 36:151 down to 37:152.
$ws.Rows("36").Insert()

# Populate the newly inserted row 36 with the new weekly data point.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44487
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100102
$ws.Range("H36").Value = "Cítricos"
$ws.Range("I36").Value = 100102006
$ws.Range("J36").Value = "Pomelo"
$ws.Range("K36").Value = "Start Ruby"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 80
$ws.Range("N36").Value = 12000
$ws.Range("O36").Value = 12000
$ws.Range("P36").Value = 12000
$ws.Range("Q36").Value = "$/bandeja 15 kilos granel"
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 800
$ws.Range("T36").Value = 15
